$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with refreshed values.
# NumberFormat is forced to text ("@") before assignment so that Excel
# keeps these numeric-looking / percent-looking strings as plain text,
# matching the inline-string cells used in the source workbook.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "329.65"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "7.07%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "39.96"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "7.00%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.259"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "2.45%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08100"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "3.78%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.529"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "2.06%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.641"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "4.73%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.918"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "1.74%"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-1.33%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9341"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "0.56%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1329"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "22.67%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1967"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "3.10%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09119"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "2.27%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03538"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "6.48%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09579"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.09%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001322"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-4.45%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.005911"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "3.61%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.367"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-5.87%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3514"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "1.15%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.980"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "11.04%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1331"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "4.42%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2560"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-1.13%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04418"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.76%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001220"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "1.25%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004319"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "1.65%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001189"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-8.70%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0003987"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-0.13%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02506"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "15.87%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05168"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "2.74%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007698"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "2.97%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1427"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "5.92%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.009191"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "5.26%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002158"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "5.71%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01099"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "37.26%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006651"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "1.27%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000749"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.45%"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "147.36%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "16.55%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002098"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.45%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0001999"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.45%"
